# Fruta / hortaliza, semanal
# Insert a new weekly observation as row 6 (pushing the existing rows 6-33
# down to 7-34) on the "Chirimoya" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row at position 6.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45230
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 21000
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 2100
$ws.Range("T6").Value = 10
